$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.463.32'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '1.618.76'
$ws.Range("E3").Value = '  +1.66%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'212.90"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = "'0.245"
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = "'0.0608"
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("D10").Value = "'19.24"
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").Value = '1.845.33'
$ws.Range("E12").Value = '  +1.58%  '

$ws.Range("D13").Value = '1.613.79'
$ws.Range("E13").Value = '  -0.25%  '

$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("D16").Value = "'63.90"
$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").Value = "'239.27"
$ws.Range("E17").Value = '  +10.92%  '

$ws.Range("D18").Value = '26.467.19'
$ws.Range("E18").Value = '  +0.81%  '

$ws.Range("D19").Value = "'7.80"
$ws.Range("E19").Value = '  +6.01%  '

$ws.Range("D20").Value = '0.0₃0726'
$ws.Range("E20").Value = '  +0.32%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("E23").Value = '  +4.28%  '

$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").Value = "'147.04"
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").Value = "'7.04"
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = '  +2.85%  '

$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").Value = '1.526.54'
$ws.Range("E32").Value = '  +7.58%  '

$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").Value = "'2.98"
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = '  +6.63%  '

$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("D37").Value = "'0.568"
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("D39").Value = "'0.831"
$ws.Range("E39").Value = '  +0.70%  '

$ws.Range("E40").Value = '  +2.46%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("D43").Value = '1.756.92'
$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("D44").Value = "'0.762"
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").Value = "'0.915"
$ws.Range("E45").Value = '  -2.42%  '

$ws.Range("D46").Value = "'61.61"
$ws.Range("E46").Value = '  +1.15%  '

$ws.Range("D47").Value = "'90.27"
$ws.Range("E47").Value = '  +3.91%  '

$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = '  +1.07%  '
